$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (ColumnWidth is stored as width-minus-5px padding,
# so subtract the standard 0.8333.. offset to land exactly on the target
# stored width)
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668
$ws.Columns.Item(3).ColumnWidth = 27.166666666666668
$ws.Columns.Item(4).ColumnWidth = 22.166666666666668
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666

# New header for column E
$ws.Range("E1").Value = "Цена за все экраны"

# Update existing row 2 values
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1000

# New row 3, only column E populated
$ws.Range("E3").Value = 1000
